# Rename the worksheet "Property1" -> "DataNode" to unify the
# conception of DataNode / DataTable / Entity (per commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Reposition the active cell/selection to C41, matching the
# author's cursor position when the workbook was last saved.
$ws.Range("C41").Select()
